$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 979
$ws.Range("J52").Value = 979
$ws.Range("L52").Value = 2937
$ws.Range("N52").Value = -3257
$ws.Range("H64").Value = 3921.5
$ws.Range("I64").Value = 3598
$ws.Range("J64").Value = 4180.3
$ws.Range("K64").Value = 3598
$ws.Range("L64").Value = 4180.3
$ws.Range("M64").Value = -3350
$ws.Range("N64").Value = -4676.3
$ws.Range("H67").Value = 3921.5
$ws.Range("I67").Value = 3598
$ws.Range("J67").Value = 4180.3
$ws.Range("K67").Value = 3598
$ws.Range("L67").Value = 4180.3
$ws.Range("M67").Value = -2740
$ws.Range("N67").Value = -5896.3
$ws.Range("H70").Value = 1946.6154
$ws.Range("J70").Value = 2118.7273
$ws.Range("L70").Value = 6356.1819
$ws.Range("N70").Value = -6896.1819
$ws.Range("H73").Value = 1946.6154
$ws.Range("J73").Value = 2118.7273
$ws.Range("L73").Value = 6356.1819
$ws.Range("N73").Value = -8228.1819
$ws.Range("H74").Value = 5520
$ws.Range("H76").Value = 4059.6667
$ws.Range("I76").Value = 3590
$ws.Range("K76").Value = 3590
$ws.Range("M76").Value = -3275
$ws.Range("H77").Value = 5520
$ws.Range("H79").Value = 4059.6667
$ws.Range("I79").Value = 3590
$ws.Range("K79").Value = 3590
$ws.Range("M79").Value = -2498
$ws.Range("H116").Value = 2132.8096
$ws.Range("I116").Value = 1908.3636
$ws.Range("J116").Value = 2379.7
$ws.Range("K116").Value = 1908.3636
$ws.Range("L116").Value = 2379.7
$ws.Range("M116").Value = 1533.6364
$ws.Range("N116").Value = -9263.700000000001
$ws.Range("H129").Value = 1046.597
$ws.Range("I129").Value = 407.25
$ws.Range("J129").Value = 1087.1904
$ws.Range("K129").Value = 1221.75
$ws.Range("L129").Value = 3261.5712
$ws.Range("M129").Value = 3778.25
$ws.Range("N129").Value = -13261.5712

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 2946
$ws.Range("I35").Value = 2946
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 2946
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = $null
$ws.Range("N35").Value = -2540
$ws.Range("H88").Value = 2000
$ws.Range("I88").Value = 2000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 2000
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = $null
$ws.Range("N88").Value = -1594
$ws.Range("H91").Value = 2000
$ws.Range("I91").Value = 2000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 2000
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = $null
$ws.Range("N91").Value = -596
$ws.Range("H122").Value = 1566.4814
$ws.Range("I122").Value = 1371.9474
$ws.Range("J122").Value = 2028.5
$ws.Range("K122").Value = 4115.8422
$ws.Range("L122").Value = 6085.5
$ws.Range("M122").Value = -1665.8422
$ws.Range("N122").Value = -10985.5
$ws.Range("H132").Value = 2153.2307
$ws.Range("I132").Value = 1460.1072
$ws.Range("J132").Value = 2961.875
$ws.Range("K132").Value = 4380.321599999999
$ws.Range("L132").Value = 8885.625
$ws.Range("M132").Value = -1850.321599999999
$ws.Range("N132").Value = -13945.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1614.8125
$ws.Range("I86").Value = 1585.3405
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 1585.3405
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -462.3405
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 1614.8125
$ws.Range("I89").Value = 1585.3405
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 7926.7025
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -2310.7025
$ws.Range("N89").Value = -26232
$ws.Range("H105").Value = 6428.5713
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 6428.5713
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = $null
$ws.Range("M105").Value = 6428.5713
$ws.Range("N105").Value = -9922.5713

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3011.3333
$ws.Range("I62").Value = 2887.75
$ws.Range("K62").Value = 2887.75
$ws.Range("M62").Value = -2263.75
$ws.Range("H65").Value = 3011.3333
$ws.Range("I65").Value = 2887.75
$ws.Range("K65").Value = 14438.75
$ws.Range("M65").Value = -11318.75
$ws.Range("H99").Value = 2222.4
$ws.Range("I99").Value = 1456
$ws.Range("J99").Value = 2733.3333
$ws.Range("K99").Value = 1456
$ws.Range("L99").Value = 2733.3333
$ws.Range("M99").Value = 42
$ws.Range("N99").Value = -5729.3333
$ws.Range("H126").Value = 2222.4
$ws.Range("I126").Value = 1456
$ws.Range("J126").Value = 2733.3333
$ws.Range("K126").Value = 4368
$ws.Range("L126").Value = 8199.999899999999
$ws.Range("M126").Value = -1898
$ws.Range("N126").Value = -13139.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 1500
$ws.Range("I59").Value = 1500
$ws.Range("K59").Value = 4500
$ws.Range("M59").Value = -3960

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5417.147
$ws.Range("I70").Value = 5002
$ws.Range("J70").Value = 5643.591
$ws.Range("K70").Value = 5002
$ws.Range("L70").Value = 5643.591
$ws.Range("M70").Value = -4732
$ws.Range("N70").Value = -6183.591
$ws.Range("H73").Value = 5417.147
$ws.Range("I73").Value = 5002
$ws.Range("J73").Value = 5643.591
$ws.Range("K73").Value = 5002
$ws.Range("L73").Value = 5643.591
$ws.Range("M73").Value = -4066
$ws.Range("N73").Value = -7515.591
$ws.Range("H80").Value = 6600.357
$ws.Range("I80").Value = 7561
$ws.Range("J80").Value = 6066.6665
$ws.Range("K80").Value = 7561
$ws.Range("L80").Value = 6066.6665
$ws.Range("M80").Value = -6563
$ws.Range("N80").Value = -8062.6665
$ws.Range("H83").Value = 6600.357
$ws.Range("I83").Value = 7561
$ws.Range("J83").Value = 6066.6665
$ws.Range("K83").Value = 37805
$ws.Range("L83").Value = 30333.3325
$ws.Range("M83").Value = -32813
$ws.Range("N83").Value = -40317.3325
$ws.Range("H132").Value = 4045.5833
$ws.Range("I132").Value = 4061
$ws.Range("J132").Value = 3999.3333
$ws.Range("K132").Value = 12183
$ws.Range("L132").Value = 11997.9999
$ws.Range("M132").Value = -9653
$ws.Range("N132").Value = -17057.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4020.5
$ws.Range("I40").Value = 3625
$ws.Range("J40").Value = 4284.1665
$ws.Range("K40").Value = 3625
$ws.Range("L40").Value = 4284.1665
$ws.Range("M40").Value = -3489
$ws.Range("N40").Value = -4556.1665
$ws.Range("H132").Value = 4483.095
$ws.Range("I132").Value = 4685.2
$ws.Range("J132").Value = 4299.364
$ws.Range("K132").Value = 14055.6
$ws.Range("L132").Value = 12898.092
$ws.Range("M132").Value = -11525.6
$ws.Range("N132").Value = -17958.092